# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 07:35"

# --- 2. Swap country labels whose rows keep their shared-string slot but now
#        point at a different country name (reordering in sharedStrings.xml) ---

# Nueva Caledonia (row 200) <-> Belice (row 201)
$ws.Range("A200").Value = "Belice"
$ws.Range("A201").Value = "Nueva Caledonia"

# Islas Virgenes Britanicas (row 213) <-> Papua Nueva Guinea (row 214)
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# Bonaire, San Eustaquio y Saba (row 215) <-> San Bartolome (row 216)
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# --- 3. Update the statistics that travel with the swapped/updated rows ---

# Row 75 - Hungria
$ws.Range("B75").Value = 3816
$ws.Range("C75").Value = 23
$ws.Range("D75").Value = 1996
$ws.Range("E75").Value = 1311
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 509

# Row 77 - Uzbekistan
$ws.Range("B77").Value = 3396
$ws.Range("C77").Value = 27
$ws.Range("E77").Value = 714

# Row 200 - now Belice
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

# Row 201 - now Nueva Caledonia
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# Row 213 - now Papua Nueva Guinea
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214 - now Islas Virgenes Britanicas
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
